$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column N: header "AQUI" (N1), matching the header style already
#     used by the rest of row 1 ---
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial() | Out-Null
$ws.Range("N1").Value = "AQUI"

# --- Row 3, col M: was a numeric time-serial (22:00:00 formatted as time);
#     becomes a plain text cell like M2. Clone M2 (already text, no custom
#     style) onto M3, then change its text. ---
$ws.Range("M2").Copy() | Out-Null
$ws.Range("M3").PasteSpecial() | Out-Null
$ws.Range("M3").Value = "22:00:00"

# --- Row 4: new data row ---
$ws.Range("A4").Value = 24.4
$ws.Range("B4").Value = 83
$ws.Range("C4").Value = 61.1
$ws.Range("D4").Value = 1026
$ws.Range("E4").Value = 106
$ws.Range("F4").Value = 163
$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 16.8
$ws.Range("I4").Value = 0.1
$ws.Range("J4").Value = 0.2

# L4 / M4 need to be plain text (not auto-parsed dates/times), so clone the
# already-text cells above them and then edit the value in place.
$ws.Range("L2").Copy() | Out-Null
$ws.Range("L4").PasteSpecial() | Out-Null

$ws.Range("M2").Copy() | Out-Null
$ws.Range("M4").PasteSpecial() | Out-Null
$ws.Range("M4").Value = "23:00:00"

$ws.Range("N4").Value = 163
